$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")
$ws.Range("D2").Value = 0.0001051705330610275
$ws.Range("E2").Value = 0.02404072415083647
$ws.Range("G2").Value = 0.001866200007498264
$ws.Range("H2").Value = 0.003437699284404516
$ws.Range("I2").Value = 0.006688360590487719
$ws.Range("J2").Value = 0.009857061319053173
$ws.Range("K2").Value = 0.0006329780444502831
$ws.Range("D3").Value = 0.001517576165497303
$ws.Range("E3").Value = 0.02896568505093455
$ws.Range("G3").Value = 0.001863546669483185
$ws.Range("H3").Value = 0.004860787652432919
$ws.Range("I3").Value = 0.007207036018371582
$ws.Range("J3").Value = 0.01275668945163488
$ws.Range("K3").Value = 0.0006483257748186588
$ws.Range("D4").Value = 0.001585181802511215
$ws.Range("E4").Value = 0.02904109004884958
$ws.Range("G4").Value = 0.001918567810207605
$ws.Range("H4").Value = 0.005087052471935749
$ws.Range("I4").Value = 0.007148382253944874
$ws.Range("J4").Value = 0.01260941568762064
$ws.Range("K4").Value = 0.0006388784386217594
$ws.Range("D5").Value = 0.0001466358080506325
$ws.Range("E5").Value = 0.02527594566345215
$ws.Range("G5").Value = 0.00186906149610877
$ws.Range("H5").Value = 0.003443620633333921
$ws.Range("I5").Value = 0.007121519185602665
$ws.Range("J5").Value = 0.010344541631639
$ws.Range("K5").Value = 0.0007301103323698044
$ws.Range("D6").Value = 0.002945591229945421
$ws.Range("E6").Value = 0.09641678910702467
$ws.Range("G6").Value = 0.003808163572102785
$ws.Range("H6").Value = 0.01098398957401514
$ws.Range("I6").Value = 0.06238697795197368
$ws.Range("J6").Value = 0.01404245849698782
$ws.Range("K6").Value = 0.00150584289804101
$ws.Range("D8").Value = 0.0001051705330610275
$ws.Range("E8").Value = 0.02404072415083647
$ws.Range("G8").Value = 0.001866200007498264
$ws.Range("H8").Value = 0.003437699284404516
$ws.Range("I8").Value = 0.006688360590487719
$ws.Range("J8").Value = 0.009857061319053173
$ws.Range("K8").Value = 0.0006329780444502831
$ws.Range("D9").Value = 0.001517576165497303
$ws.Range("E9").Value = 0.02896568505093455
$ws.Range("G9").Value = 0.001863546669483185
$ws.Range("H9").Value = 0.004860787652432919
$ws.Range("I9").Value = 0.007207036018371582
$ws.Range("J9").Value = 0.01275668945163488
$ws.Range("K9").Value = 0.0006483257748186588
$ws.Range("D10").Value = 0.001585181802511215
$ws.Range("E10").Value = 0.02904109004884958
$ws.Range("G10").Value = 0.001918567810207605
$ws.Range("H10").Value = 0.005087052471935749
$ws.Range("I10").Value = 0.007148382253944874
$ws.Range("J10").Value = 0.01260941568762064
$ws.Range("K10").Value = 0.0006388784386217594
$ws.Range("D11").Value = 0.0001466358080506325
$ws.Range("E11").Value = 0.02527594566345215
$ws.Range("G11").Value = 0.00186906149610877
$ws.Range("H11").Value = 0.003443620633333921
$ws.Range("I11").Value = 0.007121519185602665
$ws.Range("J11").Value = 0.010344541631639
$ws.Range("K11").Value = 0.0007301103323698044
$ws.Range("D12").Value = 0.002945591229945421
$ws.Range("E12").Value = 0.09641678910702467
$ws.Range("G12").Value = 0.003808163572102785
$ws.Range("H12").Value = 0.01098398957401514
$ws.Range("I12").Value = 0.06238697795197368
$ws.Range("J12").Value = 0.01404245849698782
$ws.Range("K12").Value = 0.00150584289804101
